$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the old row 5 (its data is fully superseded / no longer present).
#    Deleting it now (before any new rows are added below) cleanly drops the
#    row with no knock-on shifting since nothing exists below it yet.
# ---------------------------------------------------------------------------
$ws.Rows.Item(5).Delete()

# ---------------------------------------------------------------------------
# 2. Header row (row 1): add the three new trailing headers.
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 9).Value = "Fees paid to UKEF currency"
$ws.Cells.Item(1, 10).Value = "Payment currency"
$ws.Cells.Item(1, 11).Value = "Payment exchange rate"
$ws.Range("G1").Copy()
$ws.Range("I1:K1").PasteSpecial(-4122)
$ws.Cells.Item(1, 9).Value = "Fees paid to UKEF currency"
$ws.Cells.Item(1, 10).Value = "Payment currency"
$ws.Cells.Item(1, 11).Value = "Payment exchange rate"

# ---------------------------------------------------------------------------
# 3. Row 2 (existing "Exporter 1" row): populate new trailing columns.
# ---------------------------------------------------------------------------
$ws.Range("D2").Copy()
$ws.Range("I2:K2").PasteSpecial(-4122)
$ws.Cells.Item(2, 9).Value = "GBP"
$ws.Cells.Item(2, 10).Value = "GBP"

# ---------------------------------------------------------------------------
# 4. Row 3: replace entirely with the new "Exporter 2 GEF" data.
# ---------------------------------------------------------------------------
$ws.Cells.Item(3, 1).Value = "Exporter 2 GEF"
$ws.Cells.Item(3, 2).Value = 20001371
$ws.Cells.Item(3, 3).Value = "Exporter 2"
$ws.Range("A3").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Cells.Item(3, 3).Value = "Exporter 2"
$ws.Cells.Item(3, 4).Value = "GBP"
$ws.Cells.Item(3, 5).Value = 600000
$ws.Cells.Item(3, 6).Value = 100000
$ws.Cells.Item(3, 7).Value = 150
$ws.Cells.Item(3, 8).Value = 243
$ws.Range("D3").Copy()
$ws.Range("I3:K3").PasteSpecial(-4122)
$ws.Cells.Item(3, 9).Value = "GBP"
$ws.Cells.Item(3, 10).Value = "GBP"

# ---------------------------------------------------------------------------
# 5. Row 4: replace entirely with the new "Potato Gef" data.
# ---------------------------------------------------------------------------
$ws.Cells.Item(4, 1).Value = "Potato Gef"
$ws.Cells.Item(4, 2).Value = 20001371
$ws.Cells.Item(4, 3).Value = "Potato exporter"
$ws.Cells.Item(4, 4).Value = "GBP"
$ws.Cells.Item(4, 5).Value = 600000
$ws.Cells.Item(4, 6).Value = 100000
$ws.Cells.Item(4, 7).Value = 45
$ws.Cells.Item(4, 8).Value = 45
$ws.Range("D4").Copy()
$ws.Range("I4:K4").PasteSpecial(-4122)
$ws.Cells.Item(4, 9).Value = "EUR"
$ws.Cells.Item(4, 10).Value = "GBP"
$ws.Cells.Item(4, 11).Value = 1.17

Write-Host "edit.ps1 complete"

# ---------------------------------------------------------------------------
# 6. "Ghost" formatted-but-empty rows (6-16, 19-21) replicating the source
#    file's leftover formatting from manual spreadsheet editing.
# ---------------------------------------------------------------------------
foreach ($r in 6,9,10,11,15,16) {
    $rng = "A" + $r + ":K" + $r
    $ws.Range("A4:K4").Copy()
    $ws.Range($rng).PasteSpecial(-4122)
    $ws.Range($rng).ClearContents()
}

# Row 7: A only + I:K only (no B-H)
$ws.Range("A4").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").ClearContents()
$ws.Range("D4").Copy()
$ws.Range("I7:K7").PasteSpecial(-4122)
$ws.Range("I7:K7").ClearContents()

# Row 8: A only (bold), K only
$ws.Cells.Item(8, 1).Font.Bold = $true
$ws.Range("D4").Copy()
$ws.Range("K8").PasteSpecial(-4122)
$ws.Range("K8").ClearContents()

# Row 12: K only
$ws.Range("D4").Copy()
$ws.Range("K12").PasteSpecial(-4122)
$ws.Range("K12").ClearContents()

# Row 13: A only (bold), J only
$ws.Cells.Item(13, 1).Font.Bold = $true
$ws.Range("D4").Copy()
$ws.Range("J13").PasteSpecial(-4122)
$ws.Range("J13").ClearContents()

# Row 14: A:J only (no K)
$ws.Range("A4:J4").Copy()
$ws.Range("A14:J14").PasteSpecial(-4122)
$ws.Range("A14:J14").ClearContents()

# Row 19: B:H only (no A, no I:K)
$ws.Range("B4:H4").Copy()
$ws.Range("B19:H19").PasteSpecial(-4122)
$ws.Range("B19:H19").ClearContents()

# Row 20: A:H only (no I:K)
$ws.Range("A4:H4").Copy()
$ws.Range("A20:H20").PasteSpecial(-4122)
$ws.Range("A20:H20").ClearContents()

# Row 21: B:H only (no A, no I:K)
$ws.Range("B4:H4").Copy()
$ws.Range("B21:H21").PasteSpecial(-4122)
$ws.Range("B21:H21").ClearContents()

Write-Host "ghost rows full pattern done"
